$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ([string]::IsNullOrEmpty($val)) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Length -lt 2) { continue }
    if ($parts[0] -eq "System") { continue }

    $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
    $joined = $rotated -join ", "
    $cell.Value2 = $joined
}
